$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-decimal-looking price strings to stay as text (avoid Excel auto-numeric coercion)
$textCells = @("D4", "D5", "D6", "D14", "D19", "D20", "D21", "D23", "D24", "D25", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D43", "D44", "D47", "D49")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.228.19"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").Value = "3.056.39"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "548.73"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").Value = "139.79"
$ws.Range("E6").Value = "  +4.28%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.053.07"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "34.69"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "3.547.99"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "63.317.42"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("D17").Value = "3.057.14"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "480.81"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").Value = "13.66"
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "7.20"
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("D24").Value = "80.58"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").Value = "12.49"
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "7.90"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +4.35%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "25.94"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +5.99%  "
$ws.Range("D34").Value = "5.69"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").Value = "55.40"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").Value = "464.58"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").Value = "3.067.35"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("D44").Value = "28.28"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "116.73"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("E51").Value = "  +2.49%  "
